# sectorStocks.xlsx -- "added equally weighted index"
#
# 1. Remove the scratch "test" sheet entirely.
# 2. Trim the trailing orphan (A-only) rows from banking, RealEstate,
#    financialServices, foods-Beverages and Construction and Building.
# 3. Append the newly tracked tickers to pet-Chemicals (MFPC, MINC, EGAS,
#    ASCM) and telecom (OTMT, ETEL, GTHE).

$wb = $excel.ActiveWorkbook

# --- 1. Delete the "test" sheet -------------------------------------------
$wb.Worksheets.Item("test").Delete()

# --- 2. Drop trailing orphan rows ------------------------------------------
$wb.Worksheets.Item("banking").Rows.Item(12).Delete()
$wb.Worksheets.Item("RealEstate").Rows.Item(17).Delete()
$wb.Worksheets.Item("financialServices").Rows.Item(9).Delete()
$wb.Worksheets.Item("foods-Beverages").Rows.Item(21).Delete()
$wb.Worksheets.Item("Construction and Building").Rows.Item(17).Delete()

# --- 3. pet-Chemicals: add MFPC, MINC, EGAS, ASCM --------------------------
$pc = $wb.Worksheets.Item("pet-Chemicals")

$pc.Range("B7").Value = "MFPC"

$pc.Range("A6").Copy()
$pc.Range("A8").PasteSpecial(-4122)
$pc.Range("A8").Value = 6
$pc.Range("B8").Value = "MINC"

$pc.Range("A6").Copy()
$pc.Range("A9").PasteSpecial(-4122)
$pc.Range("A9").Value = 7
$pc.Range("B9").Value = "EGAS"

$pc.Range("A6").Copy()
$pc.Range("A10").PasteSpecial(-4122)
$pc.Range("A10").Value = 8
$pc.Range("B10").Value = "ASCM"

# --- 4. telecom: add OTMT, ETEL, GTHE --------------------------------------
$tc = $wb.Worksheets.Item("telecom")

$tc.Range("B3").Value = "OTMT"

$tc.Range("A2").Copy()
$tc.Range("A4").PasteSpecial(-4122)
$tc.Range("A4").Value = 2
$tc.Range("B4").Value = "ETEL"

$tc.Range("A2").Copy()
$tc.Range("A5").PasteSpecial(-4122)
$tc.Range("A5").Value = 3
$tc.Range("B5").Value = "GTHE"

# --- 5. Restore the original active tab (deleting the trailing "test"
#        sheet otherwise leaves the last remaining sheet selected) --------
$wb.Worksheets.Item("banking").Select()
